# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (rows 3-6), values refreshed with this week's data ---
$ws.Range("A3").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.48.1"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 97

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.1.2"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 141
$ws.Range("D4").Value = 97.09999999999999

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.60.0.10"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 244
$ws.Range("D5").Value = 97.8

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.30.0.6"
$ws.Range("B6").Value = 170
$ws.Range("C6").Value = 10884
$ws.Range("D6").Value = 98.8

# Totals row
$ws.Range("B7").Value = 175
$ws.Range("C7").Value = 11277

# --- "Good Drivers" table (rows 15-23) refreshed, dropping the 6 oldest/stale rows ---
$ws.Range("A15").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.50.1.6"
$ws.Range("B15").Value = 29259
$ws.Range("D15").Value = 99.90000000000001

$ws.Range("A16").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.8.3"
$ws.Range("B16").Value = 13028
$ws.Range("D16").Value = 100

$ws.Range("A17").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.50.2"
$ws.Range("B17").Value = 19910
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2023-11-06"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B18").Value = 10661
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2022-08-29"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B19").Value = 14239
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "'2022-05-23"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B20").Value = 265400
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "'2022-05-01"

$ws.Range("A21").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.30.1"
$ws.Range("B21").Value = 201061
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "'2020-06-01"

$ws.Range("A22").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.29.1"
$ws.Range("B22").Value = 40159
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "'2020-04-15"

$ws.Range("A23").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.14.1"
$ws.Range("B23").Value = 120862
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "'2018-05-26"

# The 6 rows that used to hold the oldest/retired driver entries (24-29) are
# gone entirely now that the window has rolled forward - remove them so the
# sheet's used range shrinks back down.
$ws.Range("A24:E29").EntireRow.Delete()
